# edit.ps1 -- reshape the single-section test document:
#   * "Test 1" -> "Test"
#   * two new blank paragraphs inserted
#   * "Test 2" -> "Sdf", flagged as a spell-check error (w:proofErr)
#   * two more new blank paragraphs inserted
#   * "Test 3" / "Test 4" collapsed into one paragraph, split into two runs
#     ("sdf" / "sdfsdfsdfsdf") in "Adobe Garamond Pro Bold", the second run
#     (and the paragraph mark) bumped to 24pt (w:sz 48 half-points)
#   * the "_GoBack" bookmark follows the text and ends up at the end of
#     that last paragraph
#
# The new paragraph/run layout mixes structural-only markup (w:proofErr,
# an empty w:pPr/w:rPr font block) that has no dedicated Word object-model
# property, so the whole body is re-expressed as WordprocessingML and
# dropped in via Range.InsertXML -- the standard COM technique for
# injecting literal OOXML -- rather than piecemeal Find/Replace calls.

$d = $word.ActiveDocument

$wNs = 'xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" ' +
       'xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" ' +
       'xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" ' +
       'xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" ' +
       'xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" ' +
       'xmlns:o="urn:schemas-microsoft-com:office:office" ' +
       'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
       'xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" ' +
       'xmlns:v="urn:schemas-microsoft-com:vml" ' +
       'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" ' +
       'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
       'xmlns:w10="urn:schemas-microsoft-com:office:word" ' +
       'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
       'xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" ' +
       'xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" ' +
       'xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" ' +
       'xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" ' +
       'xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" ' +
       'xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" ' +
       'mc:Ignorable="w14 w15 w16se wp14"'

$fontRPr = '<w:rFonts w:ascii="Adobe Garamond Pro Bold" w:hAnsi="Adobe Garamond Pro Bold"/>'

$body = '<w:p><w:r><w:t>Test</w:t></w:r></w:p>' +
        '<w:p/>' +
        '<w:p/>' +
        '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Sdf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p/>' +
        '<w:p/>' +
        '<w:p>' +
          '<w:pPr><w:rPr>' + $fontRPr + '<w:sz w:val="48"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr>' + $fontRPr + '</w:rPr><w:t>sdf</w:t></w:r>' +
          '<w:r><w:rPr>' + $fontRPr + '<w:sz w:val="48"/></w:rPr><w:t>sdfsdfsdfsdf</w:t></w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '</w:p>'

$documentXml = '<w:document ' + $wNs + '><w:body>' + $body + '</w:body></w:document>'

$pkg = '<?xml version="1.0" standalone="yes"?>' +
       '<?mso-application progid="Word.Document"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' + $documentXml + '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

# Replace the whole body content in one shot; Word keeps the trailing
# <w:sectPr> (the document only has one section) untouched.
$d.Content.InsertXML($pkg)
